# Weekly update: insert a new price record as row 64, pushing the
# existing rows 64-93 down to 65-94 (one new week of data at the top
# of this date-ordered block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 64..93 down by one row.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new week's record.
$ws.Cells.Item(64, 1).Value = 1
$ws.Cells.Item(64, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value = 45120
$ws.Cells.Item(64, 5).Value = 15
$ws.Cells.Item(64, 6).Value = 100112045
$ws.Cells.Item(64, 7).Value = "Zapallo"
$ws.Cells.Item(64, 8).Value = "Camote"
$ws.Cells.Item(64, 9).Value = "1a nueva(o)"
$ws.Cells.Item(64, 10).Value = 850
$ws.Cells.Item(64, 11).Value = 380
$ws.Cells.Item(64, 12).Value = 400
$ws.Cells.Item(64, 13).Value = 389
$ws.Cells.Item(64, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(64, 15).Value = "Perú"
$ws.Cells.Item(64, 16).Value = 389
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = "Hortaliza"
